# Auto-generated Excel COM-interop script
# Refreshes the Universalis market-price snapshot and the derived leve-profit
# figures on each job sheet of the Kujata profits workbook ("scheduled runner").
#
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 4943.0454
$ws.Range("I28").Value = 7027.467
$ws.Range("K28").Value = 7027.467
$ws.Range("M28").Value = -6542.467
# Row 103
$ws.Range("H103").Value = 1166
$ws.Range("I103").Value = 775
$ws.Range("J103").Value = 1426.6666
$ws.Range("K103").Value = 2325
$ws.Range("L103").Value = 4279.9998
$ws.Range("M103").Value = -1739
$ws.Range("N103").Value = -5451.9998
# Row 132
$ws.Range("H132").Value = 7578664
$ws.Range("I132").Value = 10419820
$ws.Range("J132").Value = 2248.1667
$ws.Range("K132").Value = 31259460
$ws.Range("L132").Value = 6744.500100000001
$ws.Range("M132").Value = -31256930
$ws.Range("N132").Value = -11804.5001

# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4634.0386
$ws.Range("I2").Value = 605.6923
$ws.Range("J2").Value = 8662.385
$ws.Range("K2").Value = 605.6923
$ws.Range("L2").Value = 8662.385
$ws.Range("M2").Value = -492.6923
$ws.Range("N2").Value = -8888.385
# Row 45
$ws.Range("H45").Value = 2302
$ws.Range("I45").Value = 2139.0908
$ws.Range("K45").Value = 2139.0908
$ws.Range("M45").Value = -1762.0908
# Row 110
$ws.Range("H110").Value = 1136.2693
$ws.Range("I110").Value = 1004.2083
$ws.Range("J110").Value = 2721
$ws.Range("K110").Value = 1004.2083
$ws.Range("L110").Value = 2721
$ws.Range("M110").Value = 1040.7917
$ws.Range("N110").Value = -6811
# Row 116
$ws.Range("H116").Value = 4634.0386
$ws.Range("I116").Value = 605.6923
$ws.Range("J116").Value = 8662.385
$ws.Range("K116").Value = 605.6923
$ws.Range("L116").Value = 8662.385
$ws.Range("M116").Value = 1688.3077
$ws.Range("N116").Value = -13250.385
# Row 122
$ws.Range("H122").Value = 759
$ws.Range("I122").Value = 703.8889
$ws.Range("K122").Value = 2111.6667
$ws.Range("M122").Value = 338.3332999999998
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4634.0386
$ws.Range("I3").Value = 605.6923
$ws.Range("J3").Value = 8662.385
$ws.Range("K3").Value = 605.6923
$ws.Range("L3").Value = 8662.385
$ws.Range("M3").Value = -491.6923
$ws.Range("N3").Value = -8890.385
# Row 94
$ws.Range("H94").Value = 8333956.5
$ws.Range("I94").Value = 11905248
$ws.Range("J94").Value = 943.3333
$ws.Range("K94").Value = 11905248
$ws.Range("L94").Value = 943.3333
$ws.Range("M94").Value = -11904797
$ws.Range("N94").Value = -1845.3333
# Row 99
$ws.Range("H99").Value = 55557076
$ws.Range("I99").Value = 55557076
$ws.Range("K99").Value = 55557076
$ws.Range("M99").Value = -55555578
# Row 105
$ws.Range("H105").Value = 90911620
$ws.Range("I105").Value = 111113600
$ws.Range("K105").Value = 111113600
$ws.Range("M105").Value = -111111853
# Row 107
$ws.Range("H107").Value = 1894.5555
$ws.Range("I107").Value = 1598.1818
$ws.Range("J107").Value = 2360.2856
$ws.Range("K107").Value = 1598.1818
$ws.Range("L107").Value = 2360.2856
$ws.Range("M107").Value = 321.8181999999999
$ws.Range("N107").Value = -6200.2856

# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
# Row 16
$ws.Range("H16").Value = 62500896
$ws.Range("I16").Value = 76923910
$ws.Range("K16").Value = 76923910
$ws.Range("M16").Value = -76923623
# Row 31
$ws.Range("H31").Value = 989.60974
$ws.Range("I31").Value = 670.7143
$ws.Range("J31").Value = 2849.8333
$ws.Range("K31").Value = 670.7143
$ws.Range("L31").Value = 2849.8333
$ws.Range("M31").Value = -375.7143
$ws.Range("N31").Value = -3439.8333
# Row 34
$ws.Range("H34").Value = 989.60974
$ws.Range("I34").Value = 670.7143
$ws.Range("J34").Value = 2849.8333
$ws.Range("K34").Value = 670.7143
$ws.Range("L34").Value = 2849.8333
$ws.Range("M34").Value = -468.7143
$ws.Range("N34").Value = -3253.8333
# Row 86
$ws.Range("H86").Value = 7432664
$ws.Range("I86").Value = 22225554
$ws.Range("J86").Value = 36219.168
$ws.Range("K86").Value = 22225554
$ws.Range("L86").Value = 36219.168
$ws.Range("M86").Value = -22224431
$ws.Range("N86").Value = -38465.168
# Row 89
$ws.Range("H89").Value = 7432664
$ws.Range("I89").Value = 22225554
$ws.Range("J89").Value = 36219.168
$ws.Range("K89").Value = 111127770
$ws.Range("L89").Value = 181095.84
$ws.Range("M89").Value = -111122154
$ws.Range("N89").Value = -192327.84
# Row 105
$ws.Range("H105").Value = 993.5
$ws.Range("I105").Value = 992.7143
$ws.Range("J105").Value = 999
$ws.Range("K105").Value = 992.7143
$ws.Range("L105").Value = 999
$ws.Range("M105").Value = 754.2857
$ws.Range("N105").Value = -4493
# Row 107
$ws.Range("H107").Value = 761.2727
$ws.Range("J107").Value = 806.3
$ws.Range("L107").Value = 806.3
$ws.Range("N107").Value = -4646.3
# Row 113
$ws.Range("H113").Value = 62500896
$ws.Range("I113").Value = 76923910
$ws.Range("K113").Value = 76923910
$ws.Range("M113").Value = -76921740
# Row 134
$ws.Range("H134").Value = 9010101
$ws.Range("I134").Value = 10753674
$ws.Range("J134").Value = 1641.6666
$ws.Range("K134").Value = 32261022
$ws.Range("L134").Value = 4924.9998
$ws.Range("M134").Value = -32258487
$ws.Range("N134").Value = -9994.9998
# Row 141
$ws.Range("H141").Value = 32597.2
$ws.Range("J141").Value = 32597.2
$ws.Range("L141").Value = 32597.2
$ws.Range("N141").Value = -42957.2

# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")
# Row 31
$ws.Range("H31").Value = 2000
$ws.Range("J31").Value = 2000
$ws.Range("L31").Value = 6000
$ws.Range("N31").Value = -6576
# Row 68
$ws.Range("H68").Value = 2078.2888
$ws.Range("J68").Value = 2151.6978
$ws.Range("L68").Value = 6455.0934
$ws.Range("N68").Value = -8077.0934
# Row 71
$ws.Range("H71").Value = 2078.2888
$ws.Range("J71").Value = 2151.6978
$ws.Range("L71").Value = 19365.2802
$ws.Range("N71").Value = -27477.2802
# Row 80
$ws.Range("H80").Value = 4491.154
$ws.Range("J80").Value = 4491.154
$ws.Range("L80").Value = 13473.462
$ws.Range("N80").Value = -15345.462
# Row 83
$ws.Range("H83").Value = 4491.154
$ws.Range("J83").Value = 4491.154
$ws.Range("L83").Value = 40420.38600000001
$ws.Range("N83").Value = -49780.38600000001
# Row 131
$ws.Range("H131").Value = 19232206
$ws.Range("J131").Value = 1576.7954
$ws.Range("L131").Value = 4730.3862
$ws.Range("N131").Value = -14810.3862

# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2420.1333
$ws.Range("J102").Value = 1617.625
$ws.Range("L102").Value = 1617.625
$ws.Range("N102").Value = -4861.625
# Row 113
$ws.Range("H113").Value = 1777.5
$ws.Range("I113").Value = 1319
$ws.Range("J113").Value = 2236
$ws.Range("K113").Value = 1319
$ws.Range("L113").Value = 2236
$ws.Range("M113").Value = 851
$ws.Range("N113").Value = -6576

# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2400
$ws.Range("I40").Value = 2333.3333
$ws.Range("K40").Value = 2333.3333
$ws.Range("M40").Value = -2197.3333
# Row 132
$ws.Range("H132").Value = 55128.527
$ws.Range("I132").Value = 2411.125
$ws.Range("J132").Value = 93468.45
$ws.Range("K132").Value = 7233.375
$ws.Range("L132").Value = 280405.35
$ws.Range("M132").Value = -4703.375
$ws.Range("N132").Value = -285465.35

# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 608.9091
$ws.Range("I107").Value = 588.3333
$ws.Range("J107").Value = 701.5
$ws.Range("K107").Value = 1764.9999
$ws.Range("L107").Value = 2104.5
$ws.Range("M107").Value = 155.0001
$ws.Range("N107").Value = -5944.5
# Row 113
$ws.Range("H113").Value = 478.46155
$ws.Range("I113").Value = 333.6
$ws.Range("J113").Value = 569
$ws.Range("K113").Value = 1000.8
$ws.Range("L113").Value = 1707
$ws.Range("M113").Value = 1169.2
$ws.Range("N113").Value = -6047
# Row 132
$ws.Range("H132").Value = 2188.1956
$ws.Range("I132").Value = 1837.1082
$ws.Range("K132").Value = 5511.3246
$ws.Range("M132").Value = -2981.3246
